$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is being restructured: two new leading id-ish columns (id /
# id_version) are pushed in front of the old "lang_code" column, several
# columns are reordered, and five brand-new audit columns are appended
# (cr_dtimes, upd_by, upd_dtimes, is_deleted, del_dtimes). Easiest/most
# reliable approach: clear the used range and rewrite both rows from
# scratch in the final layout.
# ---------------------------------------------------------------------------

$schemaJson = @'
{"$schema":"http://json-schema.org/draft-07/schema#","description":"Wuri Guinea ID Schema","additionalProperties":false,"title":"Wuri Guinea ID Schema","type":"object","definitions":{"simpleType":{"uniqueItems":true,"additionalItems":false,"type":"array","items":{"additionalProperties":false,"type":"object","required":["language","value"],"properties":{"language":{"type":"string"},"value":{"type":"string"}}}},"documentType":{"additionalProperties":false,"type":"object","properties":{"format":{"type":"string"},"type":{"type":"string"},"value":{"type":"string"},"refNumber":{"type":["string","null"]}}},"biometricsType":{"additionalProperties":false,"type":"object","properties":{"format":{"type":"string"},"version":{"type":"number","minimum":0},"value":{"type":"string"}}}},"properties":{"identity":{"additionalProperties":false,"type":"object","required":["IDSchemaVersion","firstName","lastName","dateOfBirth","gender","region","prefecture","subPrefectureOrCommune","district","sector","proofOfConsent","individualBiometrics"],"properties":{"proofOfConsent":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"gender":{"bioAttributes":[],"fieldCategory":"pvt","format":"","fieldType":"default","$ref":"#/definitions/simpleType"},"region":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"proofOfException-1":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"referenceIdentityNumber":{"bioAttributes":[],"validators":[{"validator":"^([0-9]{10,30})$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"kyc","type":"string","fieldType":"default"},"individualBiometrics":{"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/biometricsType"},"prefecture":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"subPrefectureOrCommune":{"bioAttributes":[],"validators":[{"validator":"^(?=.{0,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"district":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"sector":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"additionalAddresDetails":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"email":{"bioAttributes":[],"validators":[{"validator":"^[A-Za-z0-9_\\-]+(\\.[A-Za-z0-9_]+)*@[A-Za-z0-9_-]+(\\.[A-Za-z0-9_]+)*(\\.[a-zA-Z]{2,})$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"introducerRID":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","type":"string","fieldType":"default"},"introducerBiometrics":{"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/biometricsType"},"firstName":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"lastName":{"bioAttributes":[],"validators":[{"validator":"^(?=.{3,50}$).*","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"dateOfBirth":{"bioAttributes":[],"validators":[{"validator":"^(1869|18[7-9][0-9]|19[0-9][0-9]|20[0-9][0-9])/([0][1-9]|1[0-2])/([0][1-9]|[1-2][0-9]|3[01])$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"individualAuthBiometrics":{"bioAttributes":["leftEye","rightEye","rightIndex","rightLittle","rightRing","rightMiddle","leftIndex","leftLittle","leftRing","leftMiddle","leftThumb","rightThumb","face"],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/biometricsType"},"introducerUIN":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","type":"string","fieldType":"default"},"proofOfIdentity":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"IDSchemaVersion":{"bioAttributes":[],"fieldCategory":"none","format":"none","type":"number","fieldType":"default","minimum":0},"proofOfException":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"phone":{"bioAttributes":[],"validators":[{"validator":"^(6[1256]{1})([0-9]{7})$","arguments":[],"type":"regex"}],"fieldCategory":"pvt","format":"none","type":"string","fieldType":"default"},"introducerFirstName":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"introducerLastName":{"bioAttributes":[],"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#/definitions/simpleType"},"proofOfRelationship":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"proofOfDateOfBirth":{"bioAttributes":[],"fieldCategory":"pvt","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"proofOfAddress":{"fieldCategory":"evidence","format":"none","fieldType":"default","$ref":"#/definitions/documentType"},"UIN":{"bioAttributes":[],"fieldCategory":"none","format":"none","type":"string","fieldType":"default"}}}}}
'@

# Clear the previously used range A1:K2 before rebuilding.
$ws.Range("A1:K2").Clear()

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "id_version"
$ws.Range("C1").Value = "title"
$ws.Range("D1").Value = "description"
$ws.Range("E1").Value = "schema_json"
$ws.Range("F1").Value = "status_code"
$ws.Range("G1").Value = "add_props"
$ws.Range("H1").Value = "effective_from"
$ws.Range("I1").Value = "lang_code"
$ws.Range("J1").Value = "is_active"
$ws.Range("K1").Value = "cr_by"
$ws.Range("L1").Value = "cr_dtimes"
$ws.Range("M1").Value = "upd_by"
$ws.Range("N1").Value = "upd_dtimes"
$ws.Range("O1").Value = "is_deleted"
$ws.Range("P1").Value = "del_dtimes"

# --- Data row -----------------------------------------------------------
$ws.Range("A2").Value = 1001
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Mosip Identity"
$ws.Range("D2").Value = "Mosip Sample identity"

$ws.Range("E2").Value = $schemaJson

$ws.Range("F2").WrapText = $true
$ws.Range("F2").Value = "PUBLISHED"
$ws.Range("G2").Value = $false
$ws.Range("H2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = 0.036071759259259262
$ws.Range("I2").Value = "eng"
$ws.Range("J2").Value = $true
$ws.Range("K2").Value = "superadmin"
$ws.Range("L2").NumberFormat = "mm:ss.0"
$ws.Range("L2").Value = 0.036047453703703707
$ws.Range("M2").Value = "NULL"
$ws.Range("N2").Value = "NULL"
$ws.Range("O2").Value = $false
$ws.Range("P2").Value = "NULL"

# --- Column widths / row sizing -----------------------------------------
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Rows.Item(2).AutoFit()

# --- Selection ------------------------------------------------------------
[void]$ws.Range("E10").Select()
